$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Activecases"

$data = @(
    ,@("Aurrum Aged Care Brunswick West", 10)
    ,@("Ballarat High School Lake Gardens", 10)
    ,@("Baptcare Westhaven community", 14)
    ,@("Bialik College Hawthorn", 10)
    ,@("BlueCross Elly Kay Mordialloc", 37)
    ,@("BlueCross The Gables Camberwell", 18)
    ,@("Brookside P-9 College Caroline Springs", 20)
    ,@("Buangor Primary School", 20)
    ,@("Campbellfield Heights Primary School", 20)
    ,@("Chairo Christian School Leongatha", 11)
    ,@("Christ the Priest Primary School Caroline Springs", 18)
    ,@("Coburg North Primary School", 24)
    ,@("Covenant College Bell Post Hill", 11)
    ,@("Diamond Creek East Primary School", 12)
    ,@("Epping Views Primary School", 15)
    ,@("Estia Health Ardeer", 14)
    ,@("Fronditha Thalpori St Albans Aged Care", 20)
    ,@("Gardenvale Primary School Senior School Campus Brighton East", 16)
    ,@("Harvest Home Primary School Epping", 10)
    ,@("Hope Aged Care Sunshine West", 22)
    ,@("Ilim College Inverloch Crescent Dallas", 10)
    ,@("Kororoit Creek Primary School Burnside Heights", 19)
    ,@("Mambourin Allara Deer Park", 12)
    ,@("Marymede Catholic College South Morang", 10)
    ,@("Mornington Bay Care Community Mount Martha", 14)
    ,@("Mother of God Primary School Ardeer", 10)
    ,@("North St Kilda Children's Centre St Kilda", 19)
    ,@("Plenty Parklands Primary School Mill Park", 12)
    ,@("Rosebud Primary School", 15)
    ,@("Roxburgh Rise Primary School Roxburgh Park", 14)
    ,@("Springside Primary School Caroline Springs", 22)
    ,@("St Albans Heights Primary Schoo", 11)
    ,@("St Augustine's Primary School Yarraville", 12)
    ,@("St Brendans Primary School Shepparton", 16)
    ,@("St John's Primary School Footscray", 10)
    ,@("St Louis de Montfort's School Aspendale", 12)
    ,@("St Luke's Catholic Primary School Shepparton North", 16)
    ,@("St Mary's Primary School Hampton", 10)
    ,@("St Paul's Anglican Grammar School Warragul", 10)
    ,@("St. Francis of Assisi Primary School Mill Park", 17)
    ,@("St. Peter's Primary School Bentleigh East", 10)
    ,@("St.Thomas Aquinas College Tynong", 10)
    ,@("Templestowe Valley Primary School Templestowe Lower", 22)
    ,@("The Evelyn Hotel Fitzroy", 10)
    ,@("The Village Early Learning Centre Sandringham", 18)
    ,@("Torquay Hotel Torquay", 30)
    ,@("Village Glen Aged Care Residences Mornington", 10)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
